# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets,
# matching the upstream data refresh ("Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 6823
$ws1.Range("F11").Value = 239
$ws1.Range("F13").Value = 3149
$ws1.Range("F14").Value = 208
$ws1.Range("F15").Value = 367
$ws1.Range("F17").Value = 555
$ws1.Range("F18").Value = 17

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 6823
$ws4.Range("F15").Value = 239
$ws4.Range("F17").Value = 3149
$ws4.Range("F18").Value = 208
$ws4.Range("F19").Value = 367
$ws4.Range("F21").Value = 555
$ws4.Range("F22").Value = 17
